$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the bold/bordered header
# formatting used by the other header cells (e.g. G1) before setting its value.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data values for the new Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
